# Updated cryptos list: refresh Price (col D) and Volume(1h) (col E) values.
# NumberFormat is forced to text ("@") before writing any D value that would
# otherwise be auto-coerced to a number by Excel (e.g. "54.00" -> 54,
# "2.530" -> 2.53), so the cell keeps its original text representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.582.61"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.113.29"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +1.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.23"
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5257"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4516"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.00"
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09011"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.38"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "2.111.10"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.804"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.072"
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.67"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001176"
$ws.Range("E17").Value = "  +3.86%  "
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06707"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.35"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.322"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "30.660.49"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.78"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.392"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("D26").Value = "2.360.91"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.27"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.24"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.530"
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.96"
$ws.Range("E30").Value = "  +1.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.193"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1071"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.635"
$ws.Range("E33").Value = "  -4.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.357"
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.983"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.901"
$ws.Range("E36").Value = "  +6.09%  "
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02639"
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06837"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2326"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.64"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6863"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.267"
$ws.Range("E43").Value = "  +1.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.93"
$ws.Range("E44").Value = "  +6.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6423"
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.310"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.744"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000363"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.76"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07282"
$ws.Range("E51").Value = "  +2.63%  "
